$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D (Price) and E (Volume) columns retain text formatting so that
# numeric-looking values (e.g. "4.80", "1.00") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.726.78"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "2.455.99"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "557.68"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "161.64"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").Value = "2.454.52"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  -6.62%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -5.18%  "
$ws.Range("D13").Value = "4.80"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "2.902.16"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "68.580.74"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "0.0000168"
$ws.Range("E16").Value = "  -3.78%  "
$ws.Range("D17").Value = "23.41"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "2.459.84"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -4.84%  "
$ws.Range("D20").Value = "341.16"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "7.00"
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D22").Value = "3.79"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "6.04"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "1.85"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").Value = "66.67"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").Value = "3.65"
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("D28").Value = "2.579.33"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "8.11"
$ws.Range("E30").Value = "  -5.74%  "
$ws.Range("D31").Value = "0.0₃0815"
$ws.Range("E31").Value = "  -6.65%  "
$ws.Range("D32").Value = "7.15"
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "433.78"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("D37").Value = "157.14"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "19.04"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("D41").Value = "17.80"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "0.302"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("D43").Value = "4.41"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").Value = "1.47"
$ws.Range("E45").Value = "  -7.06%  "
$ws.Range("D46").Value = "1.10"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("D48").Value = "132.09"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("D49").Value = "3.34"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "0.481"
$ws.Range("E51").Value = "  -4.88%  "

# Restore default (General) style on the touched range so the cells remain
# styleless like the rest of the data rows (matches original workbook).
$ws.Range("D2:E51").Style = "Normal"

